$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ------------------------------------------------------------------
# Step 1: remove the old "_GoBack" bookmark that currently sits alone
# in the empty paragraph right after the (first) "con cliente pagando"
# line (that paragraph has no runs, just the bookmark pair).
# ------------------------------------------------------------------
$pagandoPara = $null
foreach ($p in $d.Paragraphs) {
    if ($pagandoPara -eq $null -and $p.Range.Text -like "*con cliente pagando*") {
        $pagandoPara = $p
    }
}

if ($pagandoPara -ne $null) {
    $bookmarkPara = $pagandoPara.Next()
    $r = $bookmarkPara.Range
    $xmlNoBookmark = "<w:p $wns><w:pPr><w:spacing w:line='240' w:lineRule='auto'/><w:rPr><w:lang w:val='es-AR'/></w:rPr></w:pPr></w:p>"
    $r.InsertXML($xmlNoBookmark)
}

# ------------------------------------------------------------------
# Step 2: insert the new "PEDIDO" paragraph, continuing the same
# numbered list (numId 2) as the paragraph that contains
# "para cada sector?? directamente para usuario...", right before the
# following empty paragraph (the one with ind left=1080).
# ------------------------------------------------------------------
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*para cada sector*directamente para usuario*") {
        $targetPara = $p
    }
}

$insPos = $targetPara.Range.Start + 1
$insRange = $d.Range($insPos, $insPos)

$newParaXml = "<w:p $wns>" +
  "<w:pPr>" +
    "<w:pStyle w:val='Prrafodelista'/>" +
    "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr>" +
    "<w:spacing w:line='240' w:lineRule='auto'/>" +
    "<w:rPr><w:lang w:val='es-AR'/></w:rPr>" +
  "</w:pPr>" +
  "<w:r><w:rPr><w:lang w:val='es-AR'/></w:rPr><w:t xml:space='preserve'>PEDIDO: que sea un </w:t></w:r>" +
  "<w:proofErr w:type='spellStart'/>" +
  "<w:r><w:rPr><w:lang w:val='es-AR'/></w:rPr><w:t>array</w:t></w:r>" +
  "<w:proofErr w:type='spellEnd'/>" +
  "<w:r><w:rPr><w:lang w:val='es-AR'/></w:rPr><w:t xml:space='preserve'> de n&#250;meros y cada </w:t></w:r>" +
  "<w:proofErr w:type='spellStart'/>" +
  "<w:r><w:rPr><w:lang w:val='es-AR'/></w:rPr><w:t>numero</w:t></w:r>" +
  "<w:proofErr w:type='spellEnd'/>" +
  "<w:r><w:rPr><w:lang w:val='es-AR'/></w:rPr><w:t xml:space='preserve'> es el ID del producto en </w:t></w:r>" +
  "<w:proofErr w:type='gramStart'/>" +
  "<w:r><w:rPr><w:lang w:val='es-AR'/></w:rPr><w:t>la BD</w:t></w:r>" +
  "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
  "<w:proofErr w:type='gramEnd'/>" +
  "</w:p>"

$insRange.InsertXML($newParaXml)
